$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.634.12"
$ws.Range("E2").Value = "  +1.81%  "
$ws.Range("D3").Value = "3.027.58"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "512.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.04%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.440"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.58"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("E10").Value = "  +1.59%  "
$ws.Range("E11").Value = "  +4.01%  "
$ws.Range("D12").Value = "3.548.64"
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("E13").Value = "  +1.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000166"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.73%  "
$ws.Range("D17").Value = "57.674.72"
$ws.Range("E17").Value = "  +1.27%  "
$ws.Range("D18").Value = "3.032.21"
$ws.Range("E18").Value = "  +1.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "332.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.53%  "
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("E23").Value = "  +4.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.82%  "
$ws.Range("E25").Value = "  +3.42%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "0.0₃0932"
$ws.Range("E27").Value = "  +1.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.45%  "
$ws.Range("E30").Value = "  +2.97%  "
$ws.Range("E31").Value = "  -1.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.75"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "155.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.90%  "
$ws.Range("E35").Value = "  +4.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.29"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "24.87"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0684"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.17%  "
$ws.Range("D39").Value = "3.066.57"
$ws.Range("E39").Value = "  +1.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.50"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.88"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.58%  "
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("D43").Value = "2.312.38"
$ws.Range("E43").Value = "  +2.44%  "
$ws.Range("E44").Value = "  +1.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.992"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0240"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.57"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.85"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0894"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.30%  "
